# Lamar added a page: insert a new "Title and Content" slide after the
# existing title slide, with the title text "Lamar added a page" and an
# otherwise-empty content placeholder.

$p = $ppt.ActivePresentation

# ppLayoutText (2) -> the "Title and Content" auto-layout, inserted as the
# second slide (index 2).
$s = $p.Slides.Add(2, 2)

# Type the title in two passes, same as a user typing "Lamar added " and
# then appending "a page" afterwards, producing two runs.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Lamar added "
[void]$s.Shapes.Item(1).TextFrame.TextRange.InsertAfter("a page")
